$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 44579
$ws.Range("K5").Value = "Modesto"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 180
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13444
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 747
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44545
$ws.Range("K6").Value = "Castle Brite"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 18000
$ws.Range("O6").Value = 19000
$ws.Range("P6").Value = 18500
$ws.Range("Q6").Value = "$/caja 15 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 1233
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("D7").Value = 44545
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 17000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 17000
$ws.Range("Q7").Value = "$/caja 15 kilos"
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1133
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44189
$ws.Range("K8").Value = "Dina"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15500
$ws.Range("Q8").Value = "$/caja 15 kilos granel"
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 1033
$ws.Range("T8").Value = 15

# Row 9
$ws.Range("D9").Value = 44189
$ws.Range("K9").Value = "Dina"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("Q9").Value = "$/caja 15 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 933
$ws.Range("T9").Value = 15
